# Update "want-to-go" (想去人数) counts in column F across all sheets
# per the refreshed stats snapshot (gh-pages rebuild at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 12927
$ws.Cells.Item(3, 6).Value = 7246
$ws.Cells.Item(10, 6).Value = 1017
$ws.Cells.Item(12, 6).Value = 363
$ws.Cells.Item(13, 6).Value = 1032
$ws.Cells.Item(16, 6).Value = 1026
$ws.Cells.Item(17, 6).Value = 512
$ws.Cells.Item(18, 6).Value = 260
$ws.Cells.Item(22, 6).Value = 316
$ws.Cells.Item(24, 6).Value = 216
$ws.Cells.Item(25, 6).Value = 384
$ws.Cells.Item(26, 6).Value = 5274
$ws.Cells.Item(27, 6).Value = 74
$ws.Cells.Item(28, 6).Value = 1457
$ws.Cells.Item(29, 6).Value = 320
$ws.Cells.Item(30, 6).Value = 1749
$ws.Cells.Item(31, 6).Value = 1749
$ws.Cells.Item(32, 6).Value = 84
$ws.Cells.Item(33, 6).Value = 72
$ws.Cells.Item(34, 6).Value = 1381
$ws.Cells.Item(35, 6).Value = 5
$ws.Cells.Item(37, 6).Value = 605
$ws.Cells.Item(39, 6).Value = 3746

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 3747
$ws.Cells.Item(4, 6).Value = 3747
$ws.Cells.Item(7, 6).Value = 67

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 9306
$ws.Cells.Item(3, 6).Value = 566
$ws.Cells.Item(4, 6).Value = 2053

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 9306
$ws.Cells.Item(3, 6).Value = 566
$ws.Cells.Item(4, 6).Value = 2053
$ws.Cells.Item(5, 6).Value = 12927
$ws.Cells.Item(6, 6).Value = 7246
$ws.Cells.Item(8, 6).Value = 3747
$ws.Cells.Item(10, 6).Value = 1017
$ws.Cells.Item(12, 6).Value = 363
$ws.Cells.Item(13, 6).Value = 1032
$ws.Cells.Item(16, 6).Value = 1026
$ws.Cells.Item(17, 6).Value = 260
$ws.Cells.Item(21, 6).Value = 316
$ws.Cells.Item(26, 6).Value = 216
$ws.Cells.Item(27, 6).Value = 384
$ws.Cells.Item(28, 6).Value = 5274
$ws.Cells.Item(29, 6).Value = 74
$ws.Cells.Item(30, 6).Value = 1457
$ws.Cells.Item(33, 6).Value = 320
$ws.Cells.Item(35, 6).Value = 1749
$ws.Cells.Item(36, 6).Value = 1749
$ws.Cells.Item(37, 6).Value = 84
$ws.Cells.Item(38, 6).Value = 72
$ws.Cells.Item(39, 6).Value = 1381
$ws.Cells.Item(40, 6).Value = 5
$ws.Cells.Item(41, 6).Value = 605
$ws.Cells.Item(48, 6).Value = 3746
